$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Greece")

$ws.Copy($null, $ws)
$nl = $wb.Worksheets.Item($wb.Worksheets.Count)
$nl.Name = "Netherlands"
$nl.Range("B4").Value = "NGC-3144/T2191"
$nl.Range("B2").Value = "Netherlands Market"

$at = $wb.Worksheets.Item($wb.Worksheets.Count)
$at.Copy($null, $at)
$at = $wb.Worksheets.Item($wb.Worksheets.Count)
$at.Name = "Austria"
$at.Range("B4").Value = "NGC-3817/T2298"
$at.Range("B2").Value = "Austria Market"

$dk = $wb.Worksheets.Item($wb.Worksheets.Count)
$dk.Copy($null, $dk)
$dk = $wb.Worksheets.Item($wb.Worksheets.Count)
$dk.Name = "Denmark"
$dk.Range("B4").Value = "NGC-2913/T2784"
$dk.Range("B2").Value = "Denmark Market"
